$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("第2次" - 2nd homework score) values keyed by row number
$d2Scores = @{
    "2" = 10
    "3" = 7
    "4" = 10
    "5" = 8
    "6" = 10
    "7" = 10
    "8" = 9
    "9" = 10
    "10" = 10
    "11" = 10
    "14" = 9
    "15" = 9
    "16" = 9.5
    "17" = 10
    "18" = 9.5
    "19" = 9.5
    "20" = 9.5
    "21" = 9.5
    "22" = 10
    "26" = 9.5
    "27" = 9.5
    "28" = 9.5
    "29" = 10
    "30" = 9
    "31" = 10
    "32" = 10
    "33" = 8
    "36" = 9.5
    "38" = 10
    "39" = 9.5
    "40" = 10
    "41" = 9.5
    "43" = 9.5
    "44" = 9.5
    "45" = 9.5
    "46" = 9.5
    "47" = 9.5
    "48" = 10
    "50" = 10
    "51" = 9
    "52" = 8.5
    "53" = 10
    "54" = 10
    "55" = 10
    "56" = 10
    "58" = 8
    "59" = 8.5
    "60" = 9
    "63" = 10
    "64" = 9
    "65" = 7
    "66" = 9
    "67" = 10
    "68" = 10
    "69" = 9
    "70" = 9.5
    "71" = 10
    "72" = 9
    "73" = 10
    "74" = 10
    "75" = 10
    "76" = 10
    "77" = 10
    "78" = 9.5
    "79" = 9
    "80" = 10
    "81" = 9.5
    "82" = 8
    "83" = 10
    "84" = 7.5
    "85" = 9
    "86" = 9.5
    "87" = 10
    "88" = 10
    "89" = 9.5
    "90" = 9
    "91" = 9.5
    "92" = 10
    "93" = 9
    "95" = 9.5
    "96" = 9
    "97" = 9
    "98" = 9.5
    "99" = 10
    "101" = 9.5
    "102" = 9.5
    "103" = 9
    "105" = 9.5
    "106" = 8
    "107" = 10
    "108" = 9.5
    "109" = 8
    "110" = 9
    "111" = 8.5
    "112" = 10
    "114" = 9.5
    "117" = 8.5
    "118" = 9
    "119" = 9
    "120" = 10
}

foreach ($row in $d2Scores.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $d2Scores[$row]
}

# Row 6 also had its column C ("第1次") score cleared in this edit
$ws.Cells.Item(6, 3).Value = $null
